# Update countries & provincias Spain
# - Swap the display order of "Santa Lucia" / "Timor Oriental" rows (201-204 block)
# - Bump the "Datos actualizados" timestamp from 06:17 to 07:34
# - Refresh the country case/death counters for India, Pakistan, Israel,
#   Kirguistan, Uzbekistan, Tailandia and Camboya

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 7 de Agosto de 2020 a las 07:34"

# --- Swap Santa Lucia / Timor Oriental rows ----------------------------
$ws.Range("A202").Value = "Timor Oriental"
$ws.Range("A203").Value = "Santa Lucia"

# --- Country figures refresh -------------------------------------------
# India (row 6)
$ws.Range("B6").Value = 2027074
$ws.Range("C6").Value = 1665
$ws.Range("D6").Value = 1378105
$ws.Range("E6").Value = 607331

# Pakistan (row 17)
$ws.Range("B17").Value = 282645
$ws.Range("C17").Value = 782
$ws.Range("D17").Value = 258099
$ws.Range("E17").Value = 18494
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = 6052

# Israel (row 36)
$ws.Range("B36").Value = 80054
$ws.Range("C36").Value = 495
$ws.Range("D36").Value = 55256
$ws.Range("E36").Value = 24222

# Kirguistan (row 56)
$ws.Range("B56").Value = 39162
$ws.Range("C56").Value = 503
$ws.Range("D56").Value = 30764
$ws.Range("E56").Value = 6947
$ws.Range("G56").Value = 4
$ws.Range("H56").Value = 1451

# Uzbekistan (row 62)
$ws.Range("B62").Value = 28535
$ws.Range("C62").Value = 220
$ws.Range("D62").Value = 19587
$ws.Range("E62").Value = 8771
$ws.Range("G62").Value = 2
$ws.Range("H62").Value = 177

# Tailandia (row 115)
$ws.Range("B115").Value = 3345
$ws.Range("C115").Value = 15
$ws.Range("E115").Value = 139

# Camboya (row 178)
$ws.Range("D178").Value = 214
$ws.Range("E178").Value = 29
